# Tutorial 6 solution update:
#  - Reformat the dates in column A from DD/MM/YYYY to DD-MM-YYYY (kept as text)
#  - Correct a few of the attendance tally columns (D..H) for the first
#    few data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new date text (slashes replaced with hyphens).
# The leading "'" forces the value to be stored as literal text instead of
# being auto-parsed into a date serial number (Excel would otherwise treat
# strings such as "01-08-2022" as a date because 01 is a valid month).
# Resetting .Style afterwards clears the quote-prefix/number-format flag
# that the text-entry leaves behind, so the cell ends up with no explicit
# style, matching a plain text cell.
$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    $cell.Value = "'" + $dates[$row]
    $cell.Style = "Normal"
}

# Updated tally values (only rows 3-6 change their D/E/G/H counts)
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0
